# Refactor eavl_fpr_latency: complete function extraction from function
#
# This script reproduces, via the Excel object model, the effect of moving
# the "fpr" column computation so it is populated like the other per-row
# metric columns, and of turning the one-off placeholder text in row 2 of
# the "sequences_results" sheet into a real numeric latency value (seconds)
# like every other row already has.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sequences_results")

# The per-row metric columns pr/rec/fpr (I/J/K) were reordered as part of
# extracting the fpr computation into its own function: the column that
# used to be "fpr" is now computed later (so it is the one that is blank),
# while pr/rec shift one column to the left to take its former place.
$ws.Cells.Item(1, 9).Value = "pr"
$ws.Cells.Item(1, 10).Value = "rec"
$ws.Cells.Item(1, 11).Value = "fpr"

# Row 2 previously held a literal placeholder string
# "29063697000 nanoseconds" in both the attack_duration (C) and
# time_to_detect (D) columns. Every other row already stores the
# duration/latency as a plain number (fraction-of-day style float)
# formatted with the "0" number format. Bring row 2 in line with that.
$latencySeconds = 0.0003363853819444445

$cC2 = $ws.Cells.Item(2, 3)
$cC2.Value = $latencySeconds
$cC2.NumberFormat = "0"

$cD2 = $ws.Cells.Item(2, 4)
$cD2.Value = $latencySeconds
$cD2.NumberFormat = "0"

# Columns (1-based): A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11 ...
# Column I was always empty (the "fpr" figure had not been wired up yet);
# it should now be populated with 0 for every data row, just like the
# neighboring pr/rec columns already are. The old column K value (which
# duplicated that same 0) is no longer produced, so it becomes blank.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = 0
    $ws.Cells.Item($r, 11).Value = ""
}

# Update the descriptive pandas-repr string stored in the summary sheet so
# it reflects that target_fpr is now a float64 series rather than a
# generic object series.
$wsSummary = $wb.Worksheets.Item("all_summary")
$cK2 = $wsSummary.Cells.Item(2, 11)
$text = $cK2.Text
$text = $text -replace "dtype: object$", "dtype: float64"
$cK2.Value = $text
